$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values to preserve exact text representation
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Cells.Item(2,4).Value = '29.213.85'
$ws.Cells.Item(2,5).Value = '  -0.61%  '
$ws.Cells.Item(3,4).Value = '1.861.50'
$ws.Cells.Item(5,4).Value = '0.7145'
$ws.Cells.Item(5,5).Value = '  -0.27%  '
$ws.Cells.Item(6,4).Value = '240.33'
$ws.Cells.Item(6,5).Value = '  +0.20%  '
$ws.Cells.Item(7,4).Value = '1.000'
$ws.Cells.Item(7,5).Value = '  +0.04%  '
$ws.Cells.Item(8,4).Value = '0.3083'
$ws.Cells.Item(9,4).Value = '0.07692'
$ws.Cells.Item(9,5).Value = '  -1.74%  '
$ws.Cells.Item(10,5).Value = '  +0.98%  '
$ws.Cells.Item(11,4).Value = '0.08325'
$ws.Cells.Item(11,5).Value = '  +0.88%  '
$ws.Cells.Item(12,4).Value = '1.905.62'
$ws.Cells.Item(12,5).Value = '  +1.72%  '
$ws.Cells.Item(13,4).Value = '0.7172'
$ws.Cells.Item(13,5).Value = '  -1.07%  '
$ws.Cells.Item(14,4).Value = '5.219'
$ws.Cells.Item(14,5).Value = '  -1.06%  '
$ws.Cells.Item(15,4).Value = '90.86'
$ws.Cells.Item(15,5).Value = '  -0.29%  '
$ws.Cells.Item(16,4).Value = '29.263.95'
$ws.Cells.Item(16,5).Value = '  -0.51%  '
$ws.Cells.Item(17,4).Value = '5.960'
$ws.Cells.Item(17,5).Value = '  +1.03%  '
$ws.Cells.Item(18,4).Value = '2.167.16'
$ws.Cells.Item(18,5).Value = '  +1.94%  '
$ws.Cells.Item(19,4).Value = '243.20'
$ws.Cells.Item(19,5).Value = '  -0.33%  '
$ws.Cells.Item(20,4).Value = '0.000007805'
$ws.Cells.Item(20,5).Value = '  -1.18%  '
$ws.Cells.Item(21,4).Value = '13.15'
$ws.Cells.Item(22,5).Value = '  +0.03%  '
$ws.Cells.Item(23,4).Value = '7.979'
$ws.Cells.Item(23,5).Value = '  +0.83%  '
$ws.Cells.Item(24,5).Value = '  +0.06%  '
$ws.Cells.Item(25,4).Value = '0.1612'
$ws.Cells.Item(25,5).Value = '  +3.46%  '
$ws.Cells.Item(26,4).Value = '162.93'
$ws.Cells.Item(26,5).Value = '  -0.48%  '
$ws.Cells.Item(27,4).Value = '8.909'
$ws.Cells.Item(27,5).Value = '  -1.09%  '
$ws.Cells.Item(28,5).Value = '  +1.51%  '
$ws.Cells.Item(29,4).Value = '1.342'
$ws.Cells.Item(29,5).Value = '  -1.19%  '
$ws.Cells.Item(30,4).Value = '4.447'
$ws.Cells.Item(30,5).Value = '  +1.64%  '
$ws.Cells.Item(31,5).Value = '  +0.41%  '
$ws.Cells.Item(32,4).Value = '4.254'
$ws.Cells.Item(32,5).Value = '  +3.23%  '
$ws.Cells.Item(33,4).Value = '0.05194'
$ws.Cells.Item(33,5).Value = '  -1.53%  '
$ws.Cells.Item(34,4).Value = '0.7930'
$ws.Cells.Item(34,5).Value = '  +10.01%  '
$ws.Cells.Item(35,4).Value = '1.928'
$ws.Cells.Item(35,5).Value = '  +0.12%  '
$ws.Cells.Item(36,4).Value = '1.171'
$ws.Cells.Item(36,5).Value = '  -2.27%  '
$ws.Cells.Item(37,4).Value = '2.686'
$ws.Cells.Item(37,5).Value = '  +0.37%  '
$ws.Cells.Item(38,4).Value = '0.01857'
$ws.Cells.Item(38,5).Value = '  -0.14%  '
$ws.Cells.Item(39,4).Value = '2.688'
$ws.Cells.Item(39,5).Value = '  -1.08%  '
$ws.Cells.Item(40,4).Value = '1.175.82'
$ws.Cells.Item(40,5).Value = '  -3.94%  '
$ws.Cells.Item(41,4).Value = '6.232'
$ws.Cells.Item(41,5).Value = '  +2.68%  '
$ws.Cells.Item(42,4).Value = '0.9016'
$ws.Cells.Item(42,5).Value = '  -0.74%  '
$ws.Cells.Item(43,4).Value = '72.92'
$ws.Cells.Item(43,5).Value = '  +0.35%  '
$ws.Cells.Item(44,4).Value = '0.9998'
$ws.Cells.Item(44,5).Value = '  -0.02%  '
$ws.Cells.Item(45,4).Value = '2.060.06'
$ws.Cells.Item(45,5).Value = '  +2.05%  '
$ws.Cells.Item(46,4).Value = '102.15'
$ws.Cells.Item(46,5).Value = '  -1.40%  '
$ws.Cells.Item(47,4).Value = '0.5198'
$ws.Cells.Item(47,5).Value = '  -2.73%  '
$ws.Cells.Item(48,4).Value = '1.776'
$ws.Cells.Item(48,5).Value = '  +0.88%  '
$ws.Cells.Item(49,4).Value = '9.386'
$ws.Cells.Item(49,5).Value = '  +1.74%  '
$ws.Cells.Item(50,2).Value = 'Aptos'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(50,4).Value = '7.077'
$ws.Cells.Item(50,5).Value = '  +0.65%  '
$ws.Cells.Item(51,2).Value = 'Frax'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(51,4).Value = '1.002'
$ws.Cells.Item(51,5).Value = '  -0.02%  '

# Revert number format/style on forced cells to keep original default styling
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
